# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Melón" / "Calameño" / "Primera"
# at row 271 (pushing the existing rows 271:376 down to 272:377) and
# populate the new row with the reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(271).Insert()

$ws.Cells.Item(271, 1).Value  = 10
$ws.Cells.Item(271, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(271, 3).Value  = "La Araucanía"
$ws.Cells.Item(271, 4).Value  = 44452
$ws.Cells.Item(271, 5).Value  = 9
$ws.Cells.Item(271, 6).Value  = 100112027
$ws.Cells.Item(271, 7).Value  = "Melón"
$ws.Cells.Item(271, 8).Value  = "Calameño"
$ws.Cells.Item(271, 9).Value  = "Primera"
$ws.Cells.Item(271, 10).Value = 500
$ws.Cells.Item(271, 11).Value = 1200
$ws.Cells.Item(271, 12).Value = 1300
$ws.Cells.Item(271, 13).Value = 1240
$ws.Cells.Item(271, 14).Value = "`$/unidad"
$ws.Cells.Item(271, 15).Value = "Brasil"
$ws.Cells.Item(271, 16).Value = 1240
$ws.Cells.Item(271, 17).Value = 1
$ws.Cells.Item(271, 18).Value = "Hortaliza"
